$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = 'Última actualización: 20:54:53'
$ws.Cells.Item(3,1).Value = 'Total filas: 505'
$ws.Cells.Item(194,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(195,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(226,1).Value = '11:56:55'
$ws.Cells.Item(226,3).Value = '10_OLMOS'
$ws.Cells.Item(226,4).Value = 91
$ws.Cells.Item(227,1).Value = '12:21:08'
$ws.Cells.Item(227,3).Value = '14_ABASTO'
$ws.Cells.Item(227,4).Value = 66
$ws.Cells.Item(242,1).Value = '12:21:08'
$ws.Cells.Item(242,3).Value = '225_GOMEZ'
$ws.Cells.Item(242,4).Value = 95
$ws.Cells.Item(243,1).Value = '12:59:47'
$ws.Cells.Item(243,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(243,4).Value = 57
$ws.Cells.Item(272,1).Value = '13:33:42'
$ws.Cells.Item(272,3).Value = '10_OLMOS'
$ws.Cells.Item(272,4).Value = 92
$ws.Cells.Item(273,1).Value = '14:45:17'
$ws.Cells.Item(273,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(273,4).Value = 20
$ws.Cells.Item(287,1).Value = '14:24:16'
$ws.Cells.Item(287,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(287,4).Value = 74
$ws.Cells.Item(288,1).Value = '14:56:20'
$ws.Cells.Item(288,3).Value = '10_OLMOS'
$ws.Cells.Item(288,4).Value = 42
$ws.Cells.Item(289,1).Value = '14:45:17'
$ws.Cells.Item(289,3).Value = '215A_EL PATO'
$ws.Cells.Item(289,4).Value = 53
$ws.Cells.Item(335,1).Value = '15:22:17'
$ws.Cells.Item(335,3).Value = '215A_EL PATO'
$ws.Cells.Item(335,4).Value = 102
$ws.Cells.Item(337,1).Value = '16:13:37'
$ws.Cells.Item(337,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(337,4).Value = 51
$ws.Cells.Item(351,1).Value = '17:14:34'
$ws.Cells.Item(351,3).Value = '10_OLMOS'
$ws.Cells.Item(351,4).Value = 20
$ws.Cells.Item(352,1).Value = '16:45:31'
$ws.Cells.Item(352,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(352,4).Value = 49
$ws.Cells.Item(364,1).Value = '17:51:34'
$ws.Cells.Item(364,3).Value = '16_SANTA ANA'
$ws.Cells.Item(364,4).Value = 1
$ws.Cells.Item(365,1).Value = '17:38:19'
$ws.Cells.Item(365,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(365,4).Value = 14
$ws.Cells.Item(366,1).Value = '15:53:28'
$ws.Cells.Item(366,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(366,4).Value = 119
$ws.Cells.Item(447,1).Value = '20:12:40'
$ws.Cells.Item(447,3).Value = '10_OLMOS'
$ws.Cells.Item(447,4).Value = 0
$ws.Cells.Item(448,1).Value = '19:13:07'
$ws.Cells.Item(448,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(448,4).Value = 59
$ws.Cells.Item(475,1).Value = '20:54:53'
$ws.Cells.Item(475,2).Value = '20:54'
$ws.Cells.Item(475,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(475,4).Value = 0
$ws.Cells.Item(476,1).Value = '20:54:53'
$ws.Cells.Item(476,2).Value = '20:54'
$ws.Cells.Item(476,3).Value = '16_SANTA ANA'
$ws.Cells.Item(476,4).Value = 0
$ws.Cells.Item(477,1).Value = '20:54:53'
$ws.Cells.Item(477,2).Value = '20:55'
$ws.Cells.Item(477,3).Value = '16_SANTA ANA'
$ws.Cells.Item(477,4).Value = 1
$ws.Cells.Item(478,1).Value = '20:12:40'
$ws.Cells.Item(478,2).Value = '20:56'
$ws.Cells.Item(478,3).Value = '27_EL RETIRO'
$ws.Cells.Item(478,4).Value = 44
$ws.Cells.Item(479,1).Value = '19:52:18'
$ws.Cells.Item(479,2).Value = '20:56'
$ws.Cells.Item(479,3).Value = '10_OLMOS'
$ws.Cells.Item(479,4).Value = 64
$ws.Cells.Item(480,1).Value = '19:13:07'
$ws.Cells.Item(480,2).Value = '20:57'
$ws.Cells.Item(480,3).Value = '27_EL RETIRO'
$ws.Cells.Item(480,4).Value = 104
$ws.Cells.Item(481,2).Value = '21:04'
$ws.Cells.Item(481,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(481,4).Value = 111
$ws.Cells.Item(482,2).Value = '21:07'
$ws.Cells.Item(482,3).Value = '215B_EL PATO'
$ws.Cells.Item(482,4).Value = 55
$ws.Cells.Item(483,1).Value = '20:33:52'
$ws.Cells.Item(483,2).Value = '21:07'
$ws.Cells.Item(483,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(483,4).Value = 34
$ws.Cells.Item(484,2).Value = '21:08'
$ws.Cells.Item(484,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(484,4).Value = 56
$ws.Cells.Item(485,1).Value = '19:13:07'
$ws.Cells.Item(485,2).Value = '21:08'
$ws.Cells.Item(485,3).Value = '215B_EL PATO'
$ws.Cells.Item(485,4).Value = 115
$ws.Cells.Item(486,1).Value = '19:39:04'
$ws.Cells.Item(486,2).Value = '21:16'
$ws.Cells.Item(486,3).Value = '14_ABASTO'
$ws.Cells.Item(486,4).Value = 97
$ws.Cells.Item(487,1).Value = '20:12:40'
$ws.Cells.Item(487,2).Value = '21:20'
$ws.Cells.Item(487,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(487,4).Value = 68
$ws.Cells.Item(488,1).Value = '19:39:04'
$ws.Cells.Item(488,2).Value = '21:21'
$ws.Cells.Item(488,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(488,4).Value = 102
$ws.Cells.Item(489,2).Value = '21:22'
$ws.Cells.Item(489,3).Value = '15_ABASTO'
$ws.Cells.Item(489,4).Value = 49
$ws.Cells.Item(490,1).Value = '19:39:04'
$ws.Cells.Item(490,2).Value = '21:23'
$ws.Cells.Item(490,3).Value = '10_OLMOS'
$ws.Cells.Item(490,4).Value = 104
$ws.Cells.Item(491,1).Value = '20:33:52'
$ws.Cells.Item(491,2).Value = '21:32'
$ws.Cells.Item(491,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(491,4).Value = 59
$ws.Cells.Item(492,1).Value = '20:33:52'
$ws.Cells.Item(492,2).Value = '21:37'
$ws.Cells.Item(492,3).Value = '14_ABASTO'
$ws.Cells.Item(492,4).Value = 64
$ws.Cells.Item(493,1).Value = '20:12:40'
$ws.Cells.Item(493,2).Value = '21:37'
$ws.Cells.Item(493,3).Value = '17_ROMERO'
$ws.Cells.Item(493,4).Value = 85
$ws.Cells.Item(494,2).Value = '21:38'
$ws.Cells.Item(494,3).Value = '14_ABASTO'
$ws.Cells.Item(494,4).Value = 106
$ws.Cells.Item(495,1).Value = '19:39:04'
$ws.Cells.Item(495,2).Value = '21:38'
$ws.Cells.Item(495,3).Value = '17_ROMERO'
$ws.Cells.Item(495,4).Value = 119
$ws.Cells.Item(496,1).Value = '20:33:52'
$ws.Cells.Item(496,2).Value = '21:46'
$ws.Cells.Item(496,3).Value = '215A_EL PATO'
$ws.Cells.Item(496,4).Value = 73
$ws.Cells.Item(497,1).Value = '19:52:18'
$ws.Cells.Item(497,2).Value = '21:47'
$ws.Cells.Item(497,3).Value = '215A_EL PATO'
$ws.Cells.Item(497,4).Value = 115
$ws.Cells.Item(498,1).Value = '20:12:40'
$ws.Cells.Item(498,2).Value = '21:52'
$ws.Cells.Item(498,3).Value = '10_OLMOS'
$ws.Cells.Item(498,4).Value = 100
$ws.Cells.Item(499,1).Value = '20:12:40'
$ws.Cells.Item(499,2).Value = '22:07'
$ws.Cells.Item(499,3).Value = '17_ROMERO'
$ws.Cells.Item(499,4).Value = 115
$ws.Cells.Item(500,2).Value = '22:07'
$ws.Cells.Item(500,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(500,4).Value = 94
$ws.Cells.Item(501,1).Value = '20:47:47'
$ws.Cells.Item(501,2).Value = '22:08'
$ws.Cells.Item(501,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(501,4).Value = 81
$ws.Cells.Item(502,2).Value = '22:13'
$ws.Cells.Item(502,3).Value = '15_ABASTO'
$ws.Cells.Item(502,4).Value = 86
$ws.Cells.Item(503,1).Value = '20:54:53'
$ws.Cells.Item(503,2).Value = '22:20'
$ws.Cells.Item(503,3).Value = '15_ABASTO'
$ws.Cells.Item(503,4).Value = 86
$ws.Cells.Item(504,1).Value = '20:33:52'
$ws.Cells.Item(504,2).Value = '22:23'
$ws.Cells.Item(504,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(504,4).Value = 110
$ws.Cells.Item(505,1).Value = '20:33:52'
$ws.Cells.Item(505,2).Value = '22:27'
$ws.Cells.Item(505,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(505,4).Value = 114
$ws.Cells.Item(505,5).Value = 'LP1912'
$ws.Cells.Item(506,1).Value = '20:47:47'
$ws.Cells.Item(506,2).Value = '22:27'
$ws.Cells.Item(506,3).Value = '10_OLMOS'
$ws.Cells.Item(506,4).Value = 100
$ws.Cells.Item(506,5).Value = 'LP1912'
$ws.Cells.Item(507,1).Value = '20:33:52'
$ws.Cells.Item(507,2).Value = '22:28'
$ws.Cells.Item(507,3).Value = '10_OLMOS'
$ws.Cells.Item(507,4).Value = 115
$ws.Cells.Item(507,5).Value = 'LP1912'
$ws.Cells.Item(508,1).Value = '20:54:53'
$ws.Cells.Item(508,2).Value = '22:38'
$ws.Cells.Item(508,3).Value = '215A_EL PATO'
$ws.Cells.Item(508,4).Value = 104
$ws.Cells.Item(508,5).Value = 'LP1912'
$ws.Cells.Item(509,1).Value = '20:47:47'
$ws.Cells.Item(509,2).Value = '22:39'
$ws.Cells.Item(509,3).Value = '215A_EL PATO'
$ws.Cells.Item(509,4).Value = 112
$ws.Cells.Item(509,5).Value = 'LP1912'
$ws.Cells.Item(510,1).Value = '20:54:53'
$ws.Cells.Item(510,2).Value = '22:47'
$ws.Cells.Item(510,3).Value = '14_ABASTO'
$ws.Cells.Item(510,4).Value = 113
$ws.Cells.Item(510,5).Value = 'LP1912'

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = 'Última actualización: 20:54:53'
$ws.Cells.Item(3,1).Value = 'Total filas: 52'
$ws.Cells.Item(56,1).Value = '20:54:53'
$ws.Cells.Item(56,2).Value = '22:38'
$ws.Cells.Item(56,4).Value = 104
$ws.Cells.Item(57,1).Value = '20:47:47'
$ws.Cells.Item(57,2).Value = '22:39'
$ws.Cells.Item(57,3).Value = '215A_EL PATO'
$ws.Cells.Item(57,4).Value = 112
$ws.Cells.Item(57,5).Value = 'LP1912'

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = 'Última actualización: 20:54:53'
